# Commit: "added css to the whole customer file"
# Net effect observed in the OOXML diff: a new order record (row 12) was
# appended to the "Order Data" sheet -- a delivery order for item [4],
# not completed, status InProgress, customer 4, waiter/chef/driver unset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Order Data")

# Columns: A Order ID | B Order Type | C Items | D Order Completed Status
#          E Order status | F Customer ID | G Waiter ID | H Chef ID | I Driver ID
$row = 12

$ws.Cells.Item($row, 1).Value = 12
$ws.Cells.Item($row, 2).Value = "delivery"
$ws.Cells.Item($row, 3).Value = "[4]"
$ws.Cells.Item($row, 4).Value = $false
$ws.Cells.Item($row, 5).Value = "InProgress"
$ws.Cells.Item($row, 6).Value = 4
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
